# Sync attendance_reports: fix "Recorded By" ordering in column G
# Rule: swap the last two comma-separated names in each cell of column G
# (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -ge 2) {
            $n = $parts.Length
            $tmp = $parts[$n - 1]
            $parts[$n - 1] = $parts[$n - 2]
            $parts[$n - 2] = $tmp

            $newVal = [string]::Join(", ", $parts)
            $cell.Value2 = $newVal
        }
    }
}
